$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.1424309830611179
    "D2" = 0.01889303280111676
    "E2" = 0.4200039547778118
    "F2" = 0.5715119509656148
    "G2" = 0.002399937969873476
    "K2" = 0.9249930234918793
    "O2" = 1.884415816162544
    "B3" = 0.1329354449600402
    "D3" = 0.0173559648507009
    "E3" = 0.3663719208139753
    "F3" = 0.5684556150569691
    "G3" = 0.002402962039085947
    "K3" = 0.8073602199524146
    "O3" = 1.889073003670148
    "B4" = 0.1271754698534835
    "D4" = 0.01640739851892903
    "E4" = 0.3335355837860021
    "F4" = 0.5671103370397077
    "G4" = 0.002404915252850254
    "K4" = 0.7348557412867649
    "O4" = 1.893720075786547
    "B5" = 0.1248460581070532
    "D5" = 0.01601967437754581
    "E5" = 0.3201762648283619
    "F5" = 0.5666952170492152
    "G5" = 0.002405735528705621
    "K5" = 0.705241098068484
    "O5" = 1.896061449082936
    "B6" = 0.1244603418153503
    "D6" = 0.01595522300728192
    "E6" = 0.3179592196631091
    "F6" = 0.5666343100681033
    "G6" = 0.002405873206197518
    "K6" = 0.7003195116135146
    "O6" = 1.896477215769778
    "B7" = 0.1271439822319138
    "D7" = 0.01640217425410384
    "E7" = 0.3333553297017033
    "F7" = 0.5671042003880302
    "G7" = 0.002404926216785786
    "K7" = 0.7344566233915941
    "O7" = 1.893749842283427
    "B8" = 0.1391424072705689
    "D8" = 0.0183640663347262
    "E8" = 0.4014909872328474
    "F8" = 0.5703475375382183
    "G8" = 0.002400960703882747
    "K8" = 0.8844915526001955
    "O8" = 1.885649498655312
    "B9" = 0.1632243973908203
    "D9" = 0.02217206104761971
    "E9" = 0.5359409132534978
    "F9" = 0.5809480528114079
    "G9" = 0.002393945807026272
    "K9" = 1.176466239640433
    "O9" = 1.884033576832849
    "B10" = 0.1812500598049525
    "D10" = 0.02494456304113157
    "E10" = 0.6353656995773207
    "F10" = 0.5913571351865841
    "G10" = 0.002389251094958128
    "K10" = 1.389578191088447
    "O10" = 1.891667137777148
    "B11" = 0.1895217846585382
    "D11" = 0.02620011994223148
    "E11" = 0.6807663500209316
    "F11" = 0.5966694471764811
    "G11" = 0.002387213961838104
    "K11" = 1.486218391214663
    "O11" = 1.897081795980483
    "B12" = 0.192664276443125
    "D12" = 0.02667472576647612
    "E12" = 0.69798535328772
    "F12" = 0.5987646778957441
    "G12" = 0.002386456636283207
    "K12" = 1.522768681585262
    "O12" = 1.89941356208746
    "B13" = 0.1919870347215209
    "D13" = 0.0265725490422497
    "E13" = 0.6942757096186085
    "F13" = 0.5983097063092515
    "G13" = 0.002386619114276583
    "K13" = 1.514898956214381
    "O13" = 1.898898828119599
    "B14" = 0.1897801164026589
    "D14" = 0.02623918318774798
    "E14" = 0.6821824195707649
    "F14" = 0.5968401448199785
    "G14" = 0.002387151374043606
    "K14" = 1.48922632010823
    "O14" = 1.897267981124287
    "B15" = 0.1884296350810075
    "D15" = 0.02603487584173791
    "E15" = 0.6747784768633522
    "F15" = 0.595950896866583
    "G15" = 0.00238747923163163
    "K15" = 1.473495156600393
    "O15" = 1.896305742479768
    "B16" = 0.1807109143098131
    "D16" = 0.02486239270292145
    "E16" = 0.6324022928535555
    "F16" = 0.5910216265279331
    "G16" = 0.002389386202445432
    "K16" = 1.38325626545452
    "O16" = 1.891352534975681
    "B17" = 0.1759940019030779
    "D17" = 0.02414163809843473
    "E17" = 0.6064513379587737
    "F17" = 0.5881459158409257
    "G17" = 0.002390581246957709
    "K17" = 1.327818432418269
    "O17" = 1.888812801517446
    "B18" = 0.1732877220365197
    "D18" = 0.02372654742808322
    "E18" = 0.5915410830662751
    "F18" = 0.5865461743056102
    "G18" = 0.002391277882063796
    "K18" = 1.295903380177492
    "O18" = 1.887534690119907
    "B19" = 0.1723725883682761
    "D19" = 0.02358591450364855
    "E19" = 0.5864954209367426
    "F19" = 0.5860138363090641
    "G19" = 0.002391515346105359
    "K19" = 1.285092607378601
    "O19" = 1.887133255144164
    "B20" = 0.1764954263836529
    "D20" = 0.02421841891092669
    "E20" = 0.6092121857102626
    "F20" = 0.5884464170364225
    "G20" = 0.002390453072823903
    "K20" = 1.333722862032459
    "O20" = 1.889064239132495
    "B21" = 0.190428067016299
    "D21" = 0.02633712403493149
    "E21" = 0.6857337691899801
    "F21" = 0.5972695178835039
    "G21" = 0.002386994654430675
    "K21" = 1.496768234329181
    "O21" = 1.897739347986118
    "B22" = 0.1995930570518851
    "D22" = 0.02771687320128535
    "E22" = 0.7359020089748896
    "F22" = 0.6035233234321566
    "G22" = 0.002384816488276058
    "K22" = 1.603063425620292
    "O22" = 1.905049941096109
    "B23" = 0.194696161680568
    "D23" = 0.02698093793772927
    "E23" = 0.7091112280774468
    "F23" = 0.6001407615152061
    "G23" = 0.002385971527721918
    "K23" = 1.546356304557946
    "O23" = 1.900997308008812
    "B24" = 0.1762687151220774
    "D24" = 0.02418370854326213
    "E24" = 0.6079639776877741
    "F24" = 0.5883103937009935
    "G24" = 0.002390510990388764
    "K24" = 1.331053603578709
    "O24" = 1.888949997350636
    "B25" = 0.1566508402819835
    "D25" = 0.02114623748008171
    "E25" = 0.499465285337223
    "F25" = 0.577622672651124
    "G25" = 0.002395762528946031
    "K25" = 1.097723150450975
    "O25" = 1.882930549535928
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
